$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Price" (D) and "Volume(1h)" (E) columns with refreshed data,
# and for three rows also update Coin (B) and Link (C) because their
# source ordering swapped with an adjacent row.
#
# Some Price values look like plain decimal numbers (e.g. "595.57") and
# Excel would otherwise auto-convert them to numeric values, stripping
# formatting such as trailing zeros. A leading apostrophe forces Excel to
# keep them as literal text, matching the original inline-string values.

$ws.Range("D2").Value = '69.560.71'
$ws.Range("E2").Value = '  +2.58%  '

$ws.Range("D3").Value = '2.502.93'
$ws.Range("E3").Value = '  +0.40%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").Value = '''595.57'
$ws.Range("E5").Value = '  +1.49%  '

$ws.Range("D6").Value = '''176.10'
$ws.Range("E6").Value = '  -0.65%  '

$ws.Range("E7").Value = '  -0.05%  '

$ws.Range("D8").Value = '''0.517'
$ws.Range("E8").Value = '  +0.47%  '

$ws.Range("D9").Value = '2.503.25'
$ws.Range("E9").Value = '  +0.44%  '

$ws.Range("D10").Value = '''0.159'
$ws.Range("E10").Value = '  +12.31%  '

$ws.Range("E11").Value = '  -0.56%  '

$ws.Range("D12").Value = '''0.342'
$ws.Range("E12").Value = '  +0.69%  '

$ws.Range("D13").Value = '''4.99'
$ws.Range("E13").Value = '  +1.26%  '

$ws.Range("D14").Value = '2.965.00'
$ws.Range("E14").Value = '  +0.63%  '

$ws.Range("D15").Value = '''25.86'
$ws.Range("E15").Value = '  +0.71%  '

$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '69.516.18'
$ws.Range("E16").Value = '  +2.62%  '

$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = '''0.0000178'
$ws.Range("E17").Value = '  +3.67%  '

$ws.Range("D18").Value = '2.499.04'
$ws.Range("E18").Value = '  +0.44%  '

$ws.Range("D19").Value = '''362.72'
$ws.Range("E19").Value = '  +3.54%  '

$ws.Range("D20").Value = '''10.99'
$ws.Range("E20").Value = '  +0.26%  '

$ws.Range("D21").Value = '''7.53'
$ws.Range("E21").Value = '  -0.18%  '

$ws.Range("D22").Value = '''4.07'
$ws.Range("E22").Value = '  -0.85%  '

$ws.Range("E23").Value = '  -0.11%  '

$ws.Range("D24").Value = '''70.64'
$ws.Range("E24").Value = '  -0.29%  '

$ws.Range("E25").Value = '  -1.81%  '

$ws.Range("D26").Value = '''9.00'
$ws.Range("E26").Value = '  -1.13%  '

$ws.Range("E27").Value = '  -3.80%  '

$ws.Range("D28").Value = '2.638.46'
$ws.Range("E28").Value = '  +0.83%  '

$ws.Range("D29").Value = '''0.997'
$ws.Range("E29").Value = '  -0.29%  '

$ws.Range("D30").Value = '''509.89'
$ws.Range("E30").Value = '  +1.07%  '

$ws.Range("D31").Value = '0.0₃0891'
$ws.Range("E31").Value = '  -1.31%  '

$ws.Range("D32").Value = '''7.73'
$ws.Range("E32").Value = '  -1.22%  '

$ws.Range("D33").Value = '''1.23'
$ws.Range("E33").Value = '  -2.19%  '

$ws.Range("E34").Value = '  +0.53%  '

$ws.Range("E35").Value = '  +0.13%  '

$ws.Range("B36").Value = 'Monero'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D36").Value = '''163.04'
$ws.Range("E36").Value = '  +0.09%  '

$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").Value = '''0.119'
$ws.Range("E37").Value = '  -2.32%  '

$ws.Range("D38").Value = '''18.74'
$ws.Range("E38").Value = '  +2.36%  '

$ws.Range("E39").Value = '  +1.19%  '

$ws.Range("B40").Value = 'USDe'
$ws.Range("C40").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D40").Value = '''1.00'
$ws.Range("E40").Value = '  +0.01%  '

$ws.Range("B41").Value = 'ImmutableX'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D41").Value = '''1.31'
$ws.Range("E41").Value = '  -1.90%  '

$ws.Range("D42").Value = '''1.72'
$ws.Range("E42").Value = '  -1.09%  '

$ws.Range("D43").Value = '''4.81'
$ws.Range("E43").Value = '  -0.97%  '

$ws.Range("D44").Value = '''0.319'
$ws.Range("E44").Value = '  -2.62%  '

$ws.Range("D45").Value = '''38.81'
$ws.Range("E45").Value = '  -0.49%  '

$ws.Range("D46").Value = '''2.32'
$ws.Range("E46").Value = '  -4.03%  '

$ws.Range("D47").Value = '''149.22'
$ws.Range("E47").Value = '  +3.18%  '

$ws.Range("D48").Value = '''3.57'
$ws.Range("E48").Value = '  +1.52%  '

$ws.Range("D49").Value = '''0.513'
$ws.Range("E49").Value = '  -0.13%  '

$ws.Range("D50").Value = '0.0₆0251'
$ws.Range("E50").Value = '  -1.11%  '

$ws.Range("D51").Value = '''0.0736'
$ws.Range("E51").Value = '  -0.61%  '
